$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7673213481903076
$ws.Range("B1").Value = 1.90214991569519
$ws.Range("C1").Value = 4.117806911468506
$ws.Range("D1").Value = 3.587493419647217
$ws.Range("E1").Value = 1.996114373207092
